# Update gh-pages to output generated at 456a3b4
#
# Applies the data refresh to the "展览" (sheet 1) and "全部类型" (sheet 4)
# worksheets: a handful of "想去人数" (F column) counters increment, and a
# brand new event row ("合肥·皖萌次元青年文化节") is appended to both sheets.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

function Set-Counter($ws, [int]$row, [double]$value) {
    $ws.Cells.Item($row, 6).Value = $value
}

function Add-NewEventRow($ws, [int]$row, [int]$indexValue) {
    $ws.Cells.Item($row, 1).Value = $indexValue

    # Column B holds a plain date-like string ("2025-02-03"), just like the
    # existing rows. Force text formatting first so Excel doesn't
    # reinterpret the literal as a real date serial, then drop back to the
    # "Normal" style so the cell ends up with no explicit style (matching
    # the sibling cells above it).
    $bCell = $ws.Cells.Item($row, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = "2025-02-03"
    $bCell.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "合肥·皖萌次元青年文化节"
    $ws.Cells.Item($row, 4).Value = "凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心"
    $ws.Cells.Item($row, 5).Value = "2025.02.03 10:00-02.04 17:30"
    $ws.Cells.Item($row, 6).Value = 2
    $ws.Cells.Item($row, 7).Value = 39.9
    $ws.Cells.Item($row, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93589"
    $ws.Cells.Item($row, 9).Value = "//i0.hdslb.com/bfs/openplatform/202410/GjWiXfOf1729133962063.jpeg"

    # Match the bold/bordered/centered format used by the other index cells
    # in column A by copying the format from the row above.
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial($xlPasteFormats)
}

# --- Sheet 1: "展览" ---------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
Set-Counter $ws1 2 441   # was 440
Set-Counter $ws1 4 3748  # was 3732
Set-Counter $ws1 5 161   # was 159
Set-Counter $ws1 6 45    # was 42
Set-Counter $ws1 7 208   # was 197
Add-NewEventRow $ws1 8 7

# --- Sheet 4: "全部类型" -----------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
Set-Counter $ws4 2 441   # was 440
Set-Counter $ws4 8 3748  # was 3732
Set-Counter $ws4 9 161   # was 159
Set-Counter $ws4 10 45   # was 42
Set-Counter $ws4 12 208  # was 197
Add-NewEventRow $ws4 13 12
